$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextValue "D2" "97.717.91"
$ws.Range("E2").Value = "  +0.89%  "
Set-TextValue "D3" "3.703.91"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +13.13%  "
Set-TextValue "D6" "237.52"
$ws.Range("E6").Value = "  -1.32%  "
Set-TextValue "D7" "656.66"
$ws.Range("E7").Value = "  -0.15%  "
Set-TextValue "D8" "0.448"
$ws.Range("E8").Value = "  +4.64%  "
$ws.Range("E9").Value = "  +3.48%  "
Set-TextValue "D10" "0.999"
$ws.Range("E10").Value = "  -0.05%  "
Set-TextValue "D11" "3.701.84"
$ws.Range("E11").Value = "  -0.53%  "
Set-TextValue "D12" "0.0000316"
$ws.Range("E12").Value = "  +16.53%  "
Set-TextValue "D13" "44.73"
$ws.Range("E13").Value = "  -1.59%  "
Set-TextValue "D14" "0.208"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("E15").Value = "  -0.86%  "
Set-TextValue "D16" "4.394.36"
Set-TextValue "D17" "97.227.42"
$ws.Range("E17").Value = "  +0.61%  "
Set-TextValue "D18" "8.89"
$ws.Range("E18").Value = "  -2.39%  "
Set-TextValue "D19" "3.708.77"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("E20").Value = "  +0.28%  "
Set-TextValue "D21" "18.77"
$ws.Range("E21").Value = "  -2.44%  "
Set-TextValue "D22" "0.536"
$ws.Range("E22").Value = "  +0.75%  "
Set-TextValue "D23" "524.82"
$ws.Range("E23").Value = "  -0.43%  "
Set-TextValue "D24" "3.44"
$ws.Range("E24").Value = "  -2.59%  "
Set-TextValue "D25" "0.0000224"
$ws.Range("E25").Value = "  +8.64%  "
Set-TextValue "D26" "118.42"
$ws.Range("E26").Value = "  +15.04%  "
Set-TextValue "D27" "6.90"
$ws.Range("E27").Value = "  -2.71%  "
Set-TextValue "D28" "0.212"
$ws.Range("E28").Value = "  +25.26%  "
Set-TextValue "D29" "13.47"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("E34").Value = "  -3.91%  "
$ws.Range("E35").Value = "  -0.33%  "
Set-TextValue "D36" "0.997"
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("E37").Value = "  -1.59%  "
Set-TextValue "D38" "637.67"
$ws.Range("E38").Value = "  -3.29%  "
Set-TextValue "D39" "8.73"
$ws.Range("E39").Value = "  -3.27%  "
Set-TextValue "D41" "0.167"
$ws.Range("E41").Value = "  +1.98%  "
Set-TextValue "D42" "0.498"
$ws.Range("E42").Value = "  +11.59%  "
$ws.Range("E43").Value = "  -5.25%  "
Set-TextValue "D44" "40.11"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("E45").Value = "  -0.32%  "
Set-TextValue "D46" "0.960"
Set-TextValue "D47" "0.0454"
$ws.Range("E47").Value = "  -1.58%  "
Set-TextValue "D48" "2.38"
$ws.Range("E48").Value = "  +0.99%  "
Set-TextValue "D49" "8.79"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("E50").Value = "  +0.09%  "
Set-TextValue "D51" "3.35"
$ws.Range("E51").Value = "  +2.84%  "
